# Add a new forecast-date column (AI, "2020-05-14") and a new observation
# row (47, "2020-05-28") to both the "cases" and "deaths" tables, and fill
# in the now-known "Observed" value for row 33 ("2020-05-14").
#
# Column layout: A = date label, B = Observed, C.. = one column per
# forecast-run date. Row layout: row 1 = headers, rows 2-46(+47) = one row
# per observation date.

$wb = $excel.ActiveWorkbook

# cases!AI34:AI46 / deaths!AI34:AI46 forecast values (new forecast-date column)
$caseForecast = @{
    34 = 19913; 35 = 20706; 36 = 21174; 37 = 21593; 38 = 22011; 39 = 22422;
    40 = 23007; 41 = 23384; 42 = 23891; 43 = 24248; 44 = 24533; 45 = 25142; 46 = 25543
}
$deathForecast = @{
    34 = 2495; 35 = 2690; 36 = 2835; 37 = 2913; 38 = 3007; 39 = 3139;
    40 = 3384; 41 = 3555; 42 = 3719; 43 = 3812; 44 = 3982; 45 = 3979; 46 = 4046
}

$sheetSpecs = @(
    @{ Name = "cases";  Forecast = $caseForecast;  ObservedB33 = 19467; NewRowValue = 25985 },
    @{ Name = "deaths"; Forecast = $deathForecast; ObservedB33 = 2247;  NewRowValue = 4150 }
)

foreach ($spec in $sheetSpecs) {
    $ws = $wb.Worksheets.Item($spec.Name)

    # New column AI (35) header: same date label as row 33's own date, "2020-05-14".
    $hdr = $ws.Cells.Item(1, 35)
    $hdr.NumberFormat = "@"
    $hdr.Value = "2020-05-14"
    $hdr.Style = "Normal"

    # Rows 2-32: new column AI has no data yet for these forecast targets.
    for ($r = 2; $r -le 32; $r++) {
        $c = $ws.Cells.Item($r, 35)
        $c.NumberFormat = "General"
        $c.Style = "Normal"
    }

    # Row 33 ("2020-05-14") now has an observed value.
    $ws.Cells.Item(33, 2).Value = $spec.ObservedB33

    # Row 33's own new forecast-column cell stays blank (forecast horizon starts the next day).
    $c33 = $ws.Cells.Item(33, 35)
    $c33.NumberFormat = "General"
    $c33.Style = "Normal"

    # Rows 34-46: fill in the new forecast column's values.
    foreach ($r in $spec.Forecast.Keys | Sort-Object) {
        $ws.Cells.Item($r, 35).Value = $spec.Forecast[$r]
    }

    # New row 47 ("2020-05-28"): only the new forecast column has a value.
    $rowLabel = $ws.Cells.Item(47, 1)
    $rowLabel.NumberFormat = "@"
    $rowLabel.Value = "2020-05-28"
    $rowLabel.Style = "Normal"

    for ($col = 2; $col -le 34; $col++) {
        $c = $ws.Cells.Item(47, $col)
        $c.NumberFormat = "General"
        $c.Style = "Normal"
    }

    $ws.Cells.Item(47, 35).Value = $spec.NewRowValue
}
